$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A4 value
$ws.Range("A4").Value = "Alberto Gomez"

# Add row 5
$ws.Range("C5:D5").NumberFormat = "@"
$ws.Range("A5").Value = "Phd. Christian Suarez"
$ws.Range("B5").Value = "Raul Alejandro Sosa"
$ws.Range("C5").Value = "172845688978"
$ws.Range("D5").Value = "1548785225"
$ws.Range("E5").Value = "Dr. Christian Santiago Izurieta Cruz"

# Add row 6
$ws.Range("C6:D6").NumberFormat = "@"
$ws.Range("A6").Value = "Lic. Pedro Peralta"
$ws.Range("B6").Value = "Josue Alberto Ramirez Arboleda"
$ws.Range("C6").Value = "174578569933"
$ws.Range("D6").Value = "1245785689"
$ws.Range("E6").Value = "Dr. Christian Santiago Izurieta Cruz"
